# Update "想去人数" (interested-count) figures to the latest scrape values.
# Sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 454
$ws1.Range("F5").Value = 216
$ws1.Range("F7").Value = 1261
$ws1.Range("F8").Value = 425
$ws1.Range("F13").Value = 437
$ws1.Range("F15").Value = 195
$ws1.Range("F27").Value = 55

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 133

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 352

# Sheet "全部类型" (All types) - aggregated view of the other three sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 352
$ws4.Range("F6").Value = 454
$ws4.Range("F7").Value = 216
$ws4.Range("F9").Value = 1261
$ws4.Range("F10").Value = 425
$ws4.Range("F20").Value = 437
$ws4.Range("F22").Value = 195
$ws4.Range("F37").Value = 133
$ws4.Range("F42").Value = 55
